$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPCOM 140")
$ws.Activate()

# Fix the opening/closing shift bug: row 25 (B25:F25) held a shift count
# of 3 but should be 2. Correcting these values also ripples through the
# SUM(B8:B40) total in B41 and the B41*5 total in B42.
$ws.Range("B25:F25").Value = 2

# Reflect the corrected row in the sheet's active selection.
$ws.Range("B25:F25").Select()
